$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.85438408664441
$ws.Range("D2").Value = 3.887263975812276
$ws.Range("E2").Value = 10.23391065261451
$ws.Range("F2").Value = 55.35832954098038
$ws.Range("G2").Value = 3.795883253417029
$ws.Range("J2").Value = 10.08419749736968
$ws.Range("K2").Value = 18.79855396859546
$ws.Range("L2").Value = 12.00897825437281
$ws.Range("M2").Value = 18.90314920690101
$ws.Range("N2").Value = 26.6832028827726

$ws.Range("B3").Value = 19.81198031517901
$ws.Range("D3").Value = 3.853773975403413
$ws.Range("E3").Value = 10.20577093064058
$ws.Range("F3").Value = 55.31422856337736
$ws.Range("G3").Value = 3.799370369121106
$ws.Range("J3").Value = 10.08003525097336
$ws.Range("K3").Value = 18.67906810664025
$ws.Range("L3").Value = 12.01341968161548
$ws.Range("M3").Value = 18.911867852338
$ws.Range("N3").Value = 26.72225449993319

$ws.Range("B4").Value = 19.79001804048887
$ws.Range("D4").Value = 3.832733214650919
$ws.Range("E4").Value = 10.18803957504637
$ws.Range("F4").Value = 55.29680060914976
$ws.Range("G4").Value = 3.801624040513216
$ws.Range("J4").Value = 10.07743125355095
$ws.Range("K4").Value = 18.61026025059196
$ws.Range("L4").Value = 12.01826256294281
$ws.Range("M4").Value = 18.9203732725609
$ws.Range("N4").Value = 26.7480985521582

$ws.Range("B5").Value = 19.78209921957695
$ws.Range("D5").Value = 3.824040357640748
$ws.Range("E5").Value = 10.18069965128882
$ws.Range("F5").Value = 55.29212839743133
$ws.Range("G5").Value = 3.802570837740805
$ws.Range("J5").Value = 10.07635780387991
$ws.Range("K5").Value = 18.58339089469751
$ws.Range("L5").Value = 12.02076850242626
$ws.Range("M5").Value = 18.92463274890514
$ws.Range("N5").Value = 26.75909968062361

$ws.Range("B6").Value = 19.78084674938633
$ws.Range("D6").Value = 3.822589799340408
$ws.Range("E6").Value = 10.17947395517894
$ws.Range("F6").Value = 55.29149940241263
$ws.Range("G6").Value = 3.802729771449201
$ws.Range("J6").Value = 10.07617881146516
$ws.Range("K6").Value = 18.5790006250752
$ws.Range("L6").Value = 12.02121677662822
$ws.Range("M6").Value = 18.92538797316971
$ws.Range("N6").Value = 26.76095477015666

$ws.Range("B7").Value = 19.78990706206961
$ws.Range("D7").Value = 3.832616457756371
$ws.Range("E7").Value = 10.18794104914608
$ws.Range("F7").Value = 55.2967277567552
$ws.Range("G7").Value = 3.801636694195207
$ws.Range("J7").Value = 10.07741682660714
$ws.Range("K7").Value = 18.60989311128595
$ws.Range("L7").Value = 12.01829420279246
$ws.Range("M7").Value = 18.92042750383051
$ws.Range("N7").Value = 26.74824501616288

$ws.Range("B8").Value = 19.83892196273257
$ws.Range("D8").Value = 3.875816714717073
$ws.Range("E8").Value = 10.22430205025382
$ws.Range("F8").Value = 55.34112291232186
$ws.Range("G8").Value = 3.797062308652345
$ws.Range("J8").Value = 10.08277230790189
$ws.Range("K8").Value = 18.75642476647455
$ws.Range("L8").Value = 12.01007079962022
$ws.Range("M8").Value = 18.90550172792485
$ws.Range("N8").Value = 26.69628066729286

$ws.Range("B9").Value = 19.96703658917564
$ws.Range("D9").Value = 3.956675176598401
$ws.Range("E9").Value = 10.29201997712414
$ws.Range("F9").Value = 55.50456332902042
$ws.Range("G9").Value = 3.788980498660356
$ws.Range("J9").Value = 10.09289722458925
$ws.Range("K9").Value = 19.07876031882079
$ws.Range("L9").Value = 12.01071000101964
$ws.Range("M9").Value = 18.90119943417015
$ws.Range("N9").Value = 26.6091777048861

$ws.Range("B10").Value = 20.08017581489867
$ws.Range("D10").Value = 4.013645486426771
$ws.Range("E10").Value = 10.3396055755231
$ws.Range("F10").Value = 55.67088383096331
$ws.Range("G10").Value = 3.783577996805101
$ws.Range("J10").Value = 10.10011558465991
$ws.Range("K10").Value = 19.33522131393641
$ws.Range("L10").Value = 12.0213605580987
$ws.Range("M10").Value = 18.91318896888932
$ws.Range("N10").Value = 26.55419374035394

$ws.Range("B11").Value = 20.13565166181533
$ws.Range("D11").Value = 4.03901612169589
$ws.Range("E11").Value = 10.36078716428695
$ws.Range("F11").Value = 55.75649105349046
$ws.Range("G11").Value = 3.781235092210371
$ws.Range("J11").Value = 10.10335361424395
$ws.Range("K11").Value = 19.45576572287248
$ws.Range("L11").Value = 12.02840411625653
$ws.Range("M11").Value = 18.92191293174158
$ws.Range("N11").Value = 26.53113449968422

$ws.Range("B12").Value = 20.1572235850172
$ws.Range("D12").Value = 4.048543357352435
$ws.Range("E12").Value = 10.36874140961438
$ws.Range("F12").Value = 55.79032787476076
$ws.Range("G12").Value = 3.78036428617902
$ws.Range("J12").Value = 10.10457336575681
$ws.Range("K12").Value = 19.501935814784
$ws.Range("L12").Value = 12.03138605474242
$ws.Range("M12").Value = 18.92568445457081
$ws.Range("N12").Value = 26.52268331385056

$ws.Range("B13").Value = 20.15255278115276
$ws.Range("D13").Value = 4.04649508375751
$ws.Range("E13").Value = 10.36703128748953
$ws.Range("F13").Value = 55.78297758222017
$ws.Range("G13").Value = 3.780551101911558
$ws.Range("J13").Value = 10.10431095370093
$ws.Range("K13").Value = 19.49196959060456
$ws.Range("L13").Value = 12.03072987010919
$ws.Range("M13").Value = 18.92485141640205
$ws.Range("N13").Value = 26.52449093743381

$ws.Range("B14").Value = 20.13741515628912
$ws.Range("D14").Value = 4.039801547174923
$ws.Range("E14").Value = 10.36144290200404
$ws.Range("F14").Value = 55.7592464547436
$ws.Range("G14").Value = 3.78116312231695
$ws.Range("J14").Value = 10.10345409178156
$ws.Range("K14").Value = 19.45955389693591
$ws.Range("L14").Value = 12.02864314422379
$ws.Range("M14").Value = 18.92221385273098
$ws.Range("N14").Value = 26.53043358704334

$ws.Range("B15").Value = 20.1282160603959
$ws.Range("D15").Value = 4.035691086708141
$ws.Range("E15").Value = 10.3580111645023
$ws.Range("F15").Value = 55.74489492335163
$ws.Range("G15").Value = 3.781540135565292
$ws.Range("J15").Value = 10.10292840728534
$ws.Range("K15").Value = 19.43976533611431
$ws.Range("L15").Value = 12.02740590035427
$ws.Range("M15").Value = 18.9206591377472
$ws.Range("N15").Value = 26.53411020561946

$ws.Range("B16").Value = 20.07662998990721
$ws.Range("D16").Value = 4.011976339234788
$ws.Range("E16").Value = 10.33821195824774
$ws.Range("F16").Value = 55.66548843780179
$ws.Range("G16").Value = 3.783733411418435
$ws.Range("J16").Value = 10.09990306563997
$ws.Range("K16").Value = 19.32741851759772
$ws.Range("L16").Value = 12.02094436608228
$ws.Range("M16").Value = 18.9126844489688
$ws.Range("N16").Value = 26.55574001778251

$ws.Range("B17").Value = 20.04600195929336
$ws.Range("D17").Value = 3.997287233150702
$ws.Range("E17").Value = 10.32594671249272
$ws.Range("F17").Value = 55.61931534459669
$ws.Range("G17").Value = 3.785108229622153
$ws.Range("J17").Value = 10.09803552883385
$ws.Range("K17").Value = 19.25946612711499
$ws.Range("L17").Value = 12.01754255957883
$ws.Range("M17").Value = 18.90862820955351
$ws.Range("N17").Value = 26.56950945784618

$ws.Range("B18").Value = 20.02876356831217
$ws.Range("D18").Value = 3.98878720436508
$ws.Range("E18").Value = 10.31884825891572
$ws.Range("F18").Value = 55.59369483797813
$ws.Range("G18").Value = 3.785909792420338
$ws.Range("J18").Value = 10.09695706677076
$ws.Range("K18").Value = 19.22074900069716
$ws.Range("L18").Value = 12.01579297012307
$ws.Range("M18").Value = 18.90660309549163
$ws.Range("N18").Value = 26.57761312488133

$ws.Range("B19").Value = 20.02299222572626
$ws.Range("D19").Value = 3.985900489983235
$ws.Range("E19").Value = 10.31643731560615
$ws.Range("F19").Value = 55.58518142856088
$ws.Range("G19").Value = 3.786183046281709
$ws.Range("J19").Value = 10.09659117444277
$ws.Range("K19").Value = 19.207704187474
$ws.Range("L19").Value = 12.01523618940432
$ws.Range("M19").Value = 18.90597037766331
$ws.Range("N19").Value = 26.58038846446917

$ws.Range("B20").Value = 20.0492233237029
$ws.Range("D20").Value = 3.998856228984828
$ws.Range("E20").Value = 10.32725690168758
$ws.Range("F20").Value = 55.62413364637437
$ws.Range("G20").Value = 3.784960760429942
$ws.Range("J20").Value = 10.09823477574339
$ws.Range("K20").Value = 19.26666201422681
$ws.Range("L20").Value = 12.01788326980569
$ws.Range("M20").Value = 18.90902814935802
$ws.Range("N20").Value = 26.56802465006587

$ws.Range("B21").Value = 20.14184622451067
$ws.Range("D21").Value = 4.041769787204713
$ws.Range("E21").Value = 10.36308615701878
$ws.Range("F21").Value = 55.76617844075506
$ws.Range("G21").Value = 3.780982912726818
$ws.Range("J21").Value = 10.10370594599606
$ws.Range("K21").Value = 19.46906127599332
$ws.Range("L21").Value = 12.02924753806087
$ws.Range("M21").Value = 18.92297588849797
$ws.Range("N21").Value = 26.52868046482043

$ws.Range("B22").Value = 20.20566445073013
$ws.Range("D22").Value = 4.06934881156766
$ws.Range("E22").Value = 10.38611356983479
$ws.Range("F22").Value = 55.86727908439365
$ws.Range("G22").Value = 3.778478713431205
$ws.Range("J22").Value = 10.10724431785265
$ws.Range("K22").Value = 19.60436968554644
$ws.Range("L22").Value = 12.03850821258394
$ws.Range("M22").Value = 18.93481785699964
$ws.Range("N22").Value = 26.50460375752517

$ws.Range("B23").Value = 20.17130714371039
$ws.Range("D23").Value = 4.054672674000048
$ws.Range("E23").Value = 10.37385895782986
$ws.Range("F23").Value = 55.812567464679
$ws.Range("G23").Value = 3.779806539958567
$ws.Range("J23").Value = 10.10535919282348
$ws.Range("K23").Value = 19.53188778728071
$ws.Range("L23").Value = 12.03339838769478
$ws.Range("M23").Value = 18.92824893388438
$ws.Range("N23").Value = 26.51730417061445

$ws.Range("B24").Value = 20.04776579122579
$ws.Range("D24").Value = 3.998147057390821
$ws.Range("E24").Value = 10.32666471158997
$ws.Range("F24").Value = 55.62195240991373
$ws.Range("G24").Value = 3.785027396477302
$ws.Range("J24").Value = 10.09814471115566
$ws.Range("K24").Value = 19.2634076620029
$ws.Range("L24").Value = 12.01772859245514
$ws.Range("M24").Value = 18.90884638053181
$ws.Range("N24").Value = 26.56869534779268

$ws.Range("B25").Value = 19.9290001097201
$ws.Range("D25").Value = 3.935222575205563
$ws.Range("E25").Value = 10.27408235412377
$ws.Range("F25").Value = 55.45219900078755
$ws.Range("G25").Value = 3.791072389041603
$ws.Range("J25").Value = 10.09019756230661
$ws.Range("K25").Value = 18.98798225824276
$ws.Range("L25").Value = 12.00874554848774
$ws.Range("M25").Value = 18.89969661259325
$ws.Range("N25").Value = 26.63115823049671
